# Daily attendance processing - 2026-01-08 04:31:11
# Swap the order of "Recorded By" values from "dnasr281@gmail.com, System"
# to "System, dnasr281@gmail.com" for every row where that exact value appears
# in column G ("Recorded By").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
